# Applies the Tab15 data refresh: updated growth-decomposition figures for
# the "Africa, Fragile States" (row 97) and "ROW, Fragile States" (row 98)
# aggregate rows, plus a mojibake fix in the Regional Economic Communities
# legend text (A103).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab15")

# --- Row 97: Africa, Fragile States ---
$ws.Range("C97").Value = 79.699764984553497
$ws.Range("D97").Value = 24.521734728529001
$ws.Range("E97").Value = -0.084002208994100006
$ws.Range("F97").Value = -4.1374975040884001
$ws.Range("G97").Value = 23.295349503536499
$ws.Range("H97").Value = 19.157851999448098

# --- Row 98: ROW, Fragile States ---
$ws.Range("C98").Value = 78.676045325895501
$ws.Range("D98").Value = 23.1393376275857
$ws.Range("E98").Value = 4.5928667145691202
$ws.Range("F98").Value = -6.4082496680503001
$ws.Range("G98").Value = 23.929108903351601
$ws.Range("H98").Value = 17.520859235301302

# --- Fix mis-encoded accented characters in the REC legend text ---
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'
